# Add a new "Serviced by " column (O) to the Card14 sheet, and fill in
# the previously-empty "Correction" column (N) with "nan" placeholders
# for the data rows, matching the rest of the sheet's convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card14")

# Header cell N1: drop the trailing space ("Correction " -> "Correction")
$ws.Range("N1").Value = "Correction"

# New header cell O1, mirroring the style of the other header cells (N1 etc.)
# Copy N1's formatting (bold font, border, center/top alignment) onto O1,
# then overwrite O1's text.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Serviced by "

# Fill N2:N12 with "nan" (same placeholder used throughout the rest of the
# table for not-yet-recorded values), and leave the new O2:O12 cells blank.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
    $ws.Cells.Item($r, 15).Value = ""
}
